# Weekly driver report refresh for 2025-04-19
#
# Table 1 ("Bad Drivers", rows 1-6) loses its third data row (the
# Intel 6E AX211 23.90.0.2 entry is gone) and gets a fresh week of
# numbers; the Totals row follows the data up.
#
# Table 2 ("Good Drivers", rows 12-22) loses its blank spacer row and
# three Wi-Fi 6E AX211 entries, keeping only six Wi-Fi 6 AX201 rows with
# refreshed sample counts / vintages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural shrink -----------------------------------------------
# Table 1: drop the old row 5 (Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.0.4)
# -- this also slides table 2 (and its leading blank rows) up by one.
$ws.Rows.Item(5).Delete()

# Table 2 data block (now rows 13-21): drop the three Wi-Fi 6E AX211
# rows so only six Wi-Fi 6 AX201 rows remain.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()

# --- column width -------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 43.17

# --- Table 1 ("Bad Drivers") refreshed data -------------------------
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.0.4"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3011
$ws.Range("D3").Value = 78.59999999999999

$ws.Range("A4").Value = "iwlwifi"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 6735
$ws.Range("D4").Value = 93.59999999999999

$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 9746

# --- Table 2 ("Good Drivers") refreshed data ------------------------
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B13").Value = 56018
$ws.Range("D13").Value = 100
$ws.Range("E13").ClearContents()

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
$ws.Range("E14").ClearContents()

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B15").Value = 442178
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").Value = "2024-11-10"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B16").Value = 77849
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "2021-08-18"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B17").Value = 59673
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "2020-08-05"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B18").Value = 113652
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "2019-12-14"
